$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '24.890.13'
Set-TextCell 'E2' '  -4.15%  '
Set-TextCell 'D3' '1.634.88'
Set-TextCell 'E3' '  -6.27%  '
Set-TextCell 'D4' '0.9969'
Set-TextCell 'E4' '  -0.23%  '
Set-TextCell 'D5' '235.37'
Set-TextCell 'E5' '  -5.22%  '
Set-TextCell 'D6' '0.9999'
Set-TextCell 'E6' '  +0.05%  '
Set-TextCell 'D7' '0.4739'
Set-TextCell 'E7' '  -6.15%  '
Set-TextCell 'D8' '0.2558'
Set-TextCell 'E8' '  -6.91%  '
Set-TextCell 'D9' '0.06063'
Set-TextCell 'E9' '  -2.00%  '
Set-TextCell 'D10' '0.06961'
Set-TextCell 'E10' '  -4.23%  '
Set-TextCell 'D11' '1.640.35'
Set-TextCell 'E11' '  -5.95%  '
Set-TextCell 'D12' '14.69'
Set-TextCell 'E12' '  -3.18%  '
Set-TextCell 'D13' '0.6103'
Set-TextCell 'E13' '  -6.67%  '
Set-TextCell 'D14' '4.338'
Set-TextCell 'E14' '  -6.65%  '
Set-TextCell 'D15' '72.89'
Set-TextCell 'E15' '  -6.14%  '
Set-TextCell 'E16' '  +0.08%  '
Set-TextCell 'D17' '0.9974'
Set-TextCell 'E17' '  -0.18%  '
Set-TextCell 'D18' '24.893.58'
Set-TextCell 'E18' '  -4.23%  '
Set-TextCell 'D19' '0.000006570'
Set-TextCell 'E19' '  -3.85%  '
Set-TextCell 'D20' '11.07'
Set-TextCell 'E20' '  -6.47%  '
Set-TextCell 'D21' '1.844.26'
Set-TextCell 'E21' '  -6.29%  '
Set-TextCell 'D22' '4.350'
Set-TextCell 'E22' '  -1.51%  '
Set-TextCell 'D23' '8.560'
Set-TextCell 'E23' '  -1.62%  '
Set-TextCell 'D24' '5.242'
Set-TextCell 'E24' '  -2.97%  '
Set-TextCell 'D25' '133.68'
Set-TextCell 'E25' '  -2.31%  '
Set-TextCell 'D26' '14.75'
Set-TextCell 'E26' '  -3.38%  '
Set-TextCell 'E27' '  -8.78%  '
Set-TextCell 'D28' '102.76'
Set-TextCell 'E28' '  -2.66%  '
Set-TextCell 'D29' '1.636'
Set-TextCell 'E29' '  -8.02%  '
Set-TextCell 'D30' '3.753'
Set-TextCell 'E30' '  -3.23%  '
Set-TextCell 'D31' '0.07716'
Set-TextCell 'E31' '  -6.24%  '
Set-TextCell 'D32' '3.539'
Set-TextCell 'E32' '  -2.70%  '
Set-TextCell 'D33' '0.9985'
Set-TextCell 'E33' '  +0.02%  '
Set-TextCell 'D34' '0.04295'
Set-TextCell 'E34' '  -8.31%  '
Set-TextCell 'D35' '2.594'
Set-TextCell 'E35' '  -2.27%  '
Set-TextCell 'D36' '0.9212'
Set-TextCell 'E36' '  -7.24%  '
Set-TextCell 'D37' '0.5793'
Set-TextCell 'E37' '  -6.45%  '
Set-TextCell 'D38' '2.557'
Set-TextCell 'E38' '  -6.96%  '
Set-TextCell 'D39' '0.01539'
Set-TextCell 'E39' '  -4.74%  '
Set-TextCell 'D40' '0.9980'
Set-TextCell 'D41' '0.8182'
Set-TextCell 'E41' '  +7.57%  '
Set-TextCell 'D42' '97.28'
Set-TextCell 'E42' '  -2.78%  '
Set-TextCell 'D43' '1.777'
Set-TextCell 'E43' '  -7.81%  '
Set-TextCell 'D44' '0.3696'
Set-TextCell 'E44' '  -6.09%  '
Set-TextCell 'D45' '4.708'
Set-TextCell 'E45' '  -6.19%  '
Set-TextCell 'B46' 'Algorand'
Set-TextCell 'C46' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D46' '0.1090'
Set-TextCell 'E46' '  -4.96%  '
Set-TextCell 'B47' 'Cronos'
Set-TextCell 'C47' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D47' '0.05195'
Set-TextCell 'E47' '  -1.48%  '
Set-TextCell 'D48' '6.021'
Set-TextCell 'E48' '  -4.55%  '
Set-TextCell 'D49' '29.44'
Set-TextCell 'E49' '  -3.92%  '
Set-TextCell 'D50' '0.9996'
Set-TextCell 'E50' '  -0.20%  '
Set-TextCell 'D51' '0.9991'
Set-TextCell 'E51' '  -0.19%  '
